$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Prepare new rows with correct header-style formatting before filling values ----
# A8 needs the same bold/border/center style as A2:A7; copy format from A7.
$ws.Cells.Item(7,1).Copy() | Out-Null
$ws.Cells.Item(8,1).PasteSpecial(-4122) | Out-Null

# J31:J37 need the same bold/border/center style as J2:J30; copy format from J30.
$ws.Cells.Item(30,10).Copy() | Out-Null
$ws.Range("J31:J37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Left table (A:H), rows 3-8 ----
$ws.Cells.Item(3,1).Value = "crude"
$ws.Cells.Item(3,2).Value = 0.8529411764705882
$ws.Cells.Item(3,3).Value = 29
$ws.Cells.Item(3,4).Value = 29
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 5

$ws.Cells.Item(4,1).Value = "crisis"
$ws.Cells.Item(4,2).Value = 0.6198630136986302
$ws.Cells.Item(4,3).Value = 181
$ws.Cells.Item(4,4).Value = 181
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 111

$ws.Cells.Item(5,1).Value = "sc"
$ws.Cells.Item(5,2).Value = 0.2328042328042328
$ws.Cells.Item(5,3).Value = 44
$ws.Cells.Item(5,4).Value = 44
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 145

$ws.Cells.Item(6,1).Value = "panic"
$ws.Cells.Item(6,2).Value = 0.2131782945736434
$ws.Cells.Item(6,3).Value = 110
$ws.Cells.Item(6,4).Value = 110
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 406

$ws.Cells.Item(7,1).Value = "low"
$ws.Cells.Item(7,2).Value = 0.1677852348993289
$ws.Cells.Item(7,3).Value = 25
$ws.Cells.Item(7,4).Value = 25
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = $false
$ws.Cells.Item(7,8).Value = 124

$ws.Cells.Item(8,1).Value = "no"
$ws.Cells.Item(8,2).Value = 0.08333333333333333
$ws.Cells.Item(8,3).Value = 30
$ws.Cells.Item(8,4).Value = 30
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = $false
$ws.Cells.Item(8,8).Value = 330

# ---- Right table (J:Q), rows 3-37 ----
$ws.Cells.Item(3,10).Value = "best"
$ws.Cells.Item(3,11).Value = 0.9491525423728814
$ws.Cells.Item(3,12).Value = 56
$ws.Cells.Item(3,13).Value = 56
$ws.Cells.Item(3,14).Value = 1
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 3

$ws.Cells.Item(4,10).Value = "love"
$ws.Cells.Item(4,11).Value = 0.8913043478260869
$ws.Cells.Item(4,12).Value = 41
$ws.Cells.Item(4,13).Value = 41
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 5

$ws.Cells.Item(5,10).Value = "special"
$ws.Cells.Item(5,11).Value = 0.8888888888888888
$ws.Cells.Item(5,12).Value = 32
$ws.Cells.Item(5,13).Value = 32
$ws.Cells.Item(5,14).Value = 1
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 4

$ws.Cells.Item(6,10).Value = "interesting"
$ws.Cells.Item(6,11).Value = 0.8787878787878788
$ws.Cells.Item(6,12).Value = 29
$ws.Cells.Item(6,13).Value = 29
$ws.Cells.Item(6,14).Value = 1
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 4

$ws.Cells.Item(7,10).Value = "great"
$ws.Cells.Item(7,11).Value = 0.8660714285714286
$ws.Cells.Item(7,12).Value = 97
$ws.Cells.Item(7,13).Value = 97
$ws.Cells.Item(7,14).Value = 1
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 15

$ws.Cells.Item(8,10).Value = "thanks"
$ws.Cells.Item(8,11).Value = 0.8292682926829268
$ws.Cells.Item(8,12).Value = 68
$ws.Cells.Item(8,13).Value = 68
$ws.Cells.Item(8,14).Value = 1
$ws.Cells.Item(8,15).Value = 0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 14

$ws.Cells.Item(9,10).Value = "thank"
$ws.Cells.Item(9,11).Value = 0.7734375
$ws.Cells.Item(9,12).Value = 99
$ws.Cells.Item(9,13).Value = 99
$ws.Cells.Item(9,14).Value = 1
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 29

$ws.Cells.Item(10,10).Value = "won"
$ws.Cells.Item(10,11).Value = 0.7692307692307693
$ws.Cells.Item(10,12).Value = 30
$ws.Cells.Item(10,13).Value = 30
$ws.Cells.Item(10,14).Value = 1
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 9

$ws.Cells.Item(11,10).Value = "safe"
$ws.Cells.Item(11,11).Value = 0.7394366197183099
$ws.Cells.Item(11,12).Value = 105
$ws.Cells.Item(11,13).Value = 105
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 37

$ws.Cells.Item(12,10).Value = "free"
$ws.Cells.Item(12,11).Value = 0.7333333333333333
$ws.Cells.Item(12,12).Value = 88
$ws.Cells.Item(12,13).Value = 88
$ws.Cells.Item(12,14).Value = 1
$ws.Cells.Item(12,15).Value = 0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 32

$ws.Cells.Item(13,10).Value = "positive"
$ws.Cells.Item(13,11).Value = 0.7241379310344828
$ws.Cells.Item(13,12).Value = 42
$ws.Cells.Item(13,13).Value = 42
$ws.Cells.Item(13,14).Value = 1
$ws.Cells.Item(13,15).Value = 0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 16

$ws.Cells.Item(14,10).Value = "good"
$ws.Cells.Item(14,11).Value = 0.7125
$ws.Cells.Item(14,12).Value = 114
$ws.Cells.Item(14,13).Value = 114
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 46

$ws.Cells.Item(15,10).Value = "confidence"
$ws.Cells.Item(15,11).Value = 0.6944444444444444
$ws.Cells.Item(15,12).Value = 25
$ws.Cells.Item(15,13).Value = 25
$ws.Cells.Item(15,14).Value = 1
$ws.Cells.Item(15,15).Value = 0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 11

$ws.Cells.Item(16,10).Value = "support"
$ws.Cells.Item(16,11).Value = 0.6792452830188679
$ws.Cells.Item(16,12).Value = 72
$ws.Cells.Item(16,13).Value = 72
$ws.Cells.Item(16,14).Value = 1
$ws.Cells.Item(16,15).Value = 0
$ws.Cells.Item(16,16).Value = $false
$ws.Cells.Item(16,17).Value = 34

$ws.Cells.Item(17,10).Value = "well"
$ws.Cells.Item(17,11).Value = 0.648936170212766
$ws.Cells.Item(17,12).Value = 61
$ws.Cells.Item(17,13).Value = 61
$ws.Cells.Item(17,14).Value = 1
$ws.Cells.Item(17,15).Value = 0
$ws.Cells.Item(17,16).Value = $false
$ws.Cells.Item(17,17).Value = 33

$ws.Cells.Item(18,10).Value = "safety"
$ws.Cells.Item(18,11).Value = 0.6078431372549019
$ws.Cells.Item(18,12).Value = 31
$ws.Cells.Item(18,13).Value = 31
$ws.Cells.Item(18,14).Value = 1
$ws.Cells.Item(18,15).Value = 0
$ws.Cells.Item(18,16).Value = $false
$ws.Cells.Item(18,17).Value = 20

$ws.Cells.Item(19,10).Value = "better"
$ws.Cells.Item(19,11).Value = 0.6031746031746031
$ws.Cells.Item(19,12).Value = 38
$ws.Cells.Item(19,13).Value = 38
$ws.Cells.Item(19,14).Value = 1
$ws.Cells.Item(19,15).Value = 0
$ws.Cells.Item(19,16).Value = $false
$ws.Cells.Item(19,17).Value = 25

$ws.Cells.Item(20,10).Value = "heroes"
$ws.Cells.Item(20,11).Value = 0.5957446808510638
$ws.Cells.Item(20,12).Value = 28
$ws.Cells.Item(20,13).Value = 28
$ws.Cells.Item(20,14).Value = 1
$ws.Cells.Item(20,15).Value = 0
$ws.Cells.Item(20,16).Value = $false
$ws.Cells.Item(20,17).Value = 19

$ws.Cells.Item(21,10).Value = "hand"
$ws.Cells.Item(21,11).Value = 0.5535248041775457
$ws.Cells.Item(21,12).Value = 212
$ws.Cells.Item(21,13).Value = 212
$ws.Cells.Item(21,14).Value = 1
$ws.Cells.Item(21,15).Value = 0
$ws.Cells.Item(21,16).Value = $false
$ws.Cells.Item(21,17).Value = 171

$ws.Cells.Item(22,10).Value = "care"
$ws.Cells.Item(22,11).Value = 0.4719101123595505
$ws.Cells.Item(22,12).Value = 42
$ws.Cells.Item(22,13).Value = 42
$ws.Cells.Item(22,14).Value = 1
$ws.Cells.Item(22,15).Value = 0
$ws.Cells.Item(22,16).Value = $false
$ws.Cells.Item(22,17).Value = 47

$ws.Cells.Item(23,10).Value = "like"
$ws.Cells.Item(23,11).Value = 0.4676470588235294
$ws.Cells.Item(23,12).Value = 159
$ws.Cells.Item(23,13).Value = 159
$ws.Cells.Item(23,14).Value = 1
$ws.Cells.Item(23,15).Value = 0
$ws.Cells.Item(23,16).Value = $false
$ws.Cells.Item(23,17).Value = 181

$ws.Cells.Item(24,10).Value = "help"
$ws.Cells.Item(24,11).Value = 0.4406779661016949
$ws.Cells.Item(24,12).Value = 130
$ws.Cells.Item(24,13).Value = 130
$ws.Cells.Item(24,14).Value = 1
$ws.Cells.Item(24,15).Value = 0
$ws.Cells.Item(24,16).Value = $false
$ws.Cells.Item(24,17).Value = 165

$ws.Cells.Item(25,10).Value = "increase"
$ws.Cells.Item(25,11).Value = 0.4230769230769231
$ws.Cells.Item(25,12).Value = 33
$ws.Cells.Item(25,13).Value = 33
$ws.Cells.Item(25,14).Value = 1
$ws.Cells.Item(25,15).Value = 0
$ws.Cells.Item(25,16).Value = $false
$ws.Cells.Item(25,17).Value = 45

$ws.Cells.Item(26,10).Value = "protect"
$ws.Cells.Item(26,11).Value = 0.410958904109589
$ws.Cells.Item(26,12).Value = 30
$ws.Cells.Item(26,13).Value = 30
$ws.Cells.Item(26,14).Value = 1
$ws.Cells.Item(26,15).Value = 0
$ws.Cells.Item(26,16).Value = $false
$ws.Cells.Item(26,17).Value = 43

$ws.Cells.Item(27,10).Value = "please"
$ws.Cells.Item(27,11).Value = 0.3891213389121339
$ws.Cells.Item(27,12).Value = 93
$ws.Cells.Item(27,13).Value = 93
$ws.Cells.Item(27,14).Value = 1
$ws.Cells.Item(27,15).Value = 0
$ws.Cells.Item(27,16).Value = $false
$ws.Cells.Item(27,17).Value = 146

$ws.Cells.Item(28,10).Value = "store"
$ws.Cells.Item(28,11).Value = 0.04026845637583892
$ws.Cells.Item(28,12).Value = 36
$ws.Cells.Item(28,13).Value = 36
$ws.Cells.Item(28,14).Value = 1
$ws.Cells.Item(28,15).Value = 0
$ws.Cells.Item(28,16).Value = $false
$ws.Cells.Item(28,17).Value = 858

$ws.Cells.Item(29,10).Value = "you"
$ws.Cells.Item(29,11).Value = 0.03583333333333334
$ws.Cells.Item(29,12).Value = 43
$ws.Cells.Item(29,13).Value = 43
$ws.Cells.Item(29,14).Value = 1
$ws.Cells.Item(29,15).Value = 0
$ws.Cells.Item(29,16).Value = $false
$ws.Cells.Item(29,17).Value = 1157

$ws.Cells.Item(30,10).Value = "grocery"
$ws.Cells.Item(30,11).Value = 0.02774694783573807
$ws.Cells.Item(30,12).Value = 25
$ws.Cells.Item(30,13).Value = 25
$ws.Cells.Item(30,14).Value = 1
$ws.Cells.Item(30,15).Value = 0
$ws.Cells.Item(30,16).Value = $false
$ws.Cells.Item(30,17).Value = 876

$ws.Cells.Item(31,10).Value = "and"
$ws.Cells.Item(31,11).Value = 0.01871958068139274
$ws.Cells.Item(31,12).Value = 50
$ws.Cells.Item(31,13).Value = 52
$ws.Cells.Item(31,14).Value = 0.96
$ws.Cells.Item(31,15).Value = 0.04000000000000004
$ws.Cells.Item(31,16).Value = $true
$ws.Cells.Item(31,17).Value = 2621

$ws.Cells.Item(32,10).Value = ","
$ws.Cells.Item(32,11).Value = 0.01351904956984842
$ws.Cells.Item(32,12).Value = 33
$ws.Cells.Item(32,13).Value = 37
$ws.Cells.Item(32,14).Value = 0.89
$ws.Cells.Item(32,15).Value = 0.11
$ws.Cells.Item(32,16).Value = $true
$ws.Cells.Item(32,17).Value = 2408

$ws.Cells.Item(33,10).Value = "to"
$ws.Cells.Item(33,11).Value = 0.01295696436834799
$ws.Cells.Item(33,12).Value = 56
$ws.Cells.Item(33,13).Value = 62
$ws.Cells.Item(33,14).Value = 0.9
$ws.Cells.Item(33,15).Value = 0.09999999999999998
$ws.Cells.Item(33,16).Value = $true
$ws.Cells.Item(33,17).Value = 4266

$ws.Cells.Item(34,10).Value = "19"
$ws.Cells.Item(34,11).Value = 0.01260504201680672
$ws.Cells.Item(34,12).Value = 27
$ws.Cells.Item(34,13).Value = 28
$ws.Cells.Item(34,14).Value = 0.96
$ws.Cells.Item(34,15).Value = 0.04000000000000004
$ws.Cells.Item(34,16).Value = $true
$ws.Cells.Item(34,17).Value = 2115

$ws.Cells.Item(35,10).Value = "."
$ws.Cells.Item(35,11).Value = 0.01141369643572287
$ws.Cells.Item(35,12).Value = 57
$ws.Cells.Item(35,13).Value = 59
$ws.Cells.Item(35,14).Value = 0.97
$ws.Cells.Item(35,15).Value = 0.03000000000000003
$ws.Cells.Item(35,16).Value = $true
$ws.Cells.Item(35,17).Value = 4937

$ws.Cells.Item(36,10).Value = "co"
$ws.Cells.Item(36,11).Value = 0.009029345372460496
$ws.Cells.Item(36,12).Value = 28
$ws.Cells.Item(36,13).Value = 34
$ws.Cells.Item(36,14).Value = 0.82
$ws.Cells.Item(36,15).Value = 0.18
$ws.Cells.Item(36,16).Value = $true
$ws.Cells.Item(36,17).Value = 3073

$ws.Cells.Item(37,10).Value = "the"
$ws.Cells.Item(37,11).Value = 0.008137957760124007
$ws.Cells.Item(37,12).Value = 42
$ws.Cells.Item(37,13).Value = 46
$ws.Cells.Item(37,14).Value = 0.91
$ws.Cells.Item(37,15).Value = 0.08999999999999997
$ws.Cells.Item(37,16).Value = $true
$ws.Cells.Item(37,17).Value = 5119
